$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.315.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.40%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.840.38'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.85%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.62%  '

$ws.Range('E6').Value = '  -1.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07406'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.57%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2895'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.24%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.79'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.59%  '

$ws.Range('E11').Value = '  -0.31%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.843.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.68%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.972'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6770'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001016'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.19%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.73%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.233'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.00%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.334.79'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.32%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.94%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9995'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.409'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.464'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.24%  '

$ws.Range('E26').Value = '  -3.66%  '

$ws.Range('E27').Value = '  -1.87%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.447'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.47%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06458'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +13.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.483'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.066'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.062'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.14%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.837'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.11%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.137'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.83%  '

$ws.Range('E35').Value = '  -1.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.567'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.74%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01857'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.240.36'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.01%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.732'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.56%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9311'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.75%  '

$ws.Range('B42').Value = 'RocketPoolETH'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.027.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.54%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9991'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.13%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.47%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.56'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000118'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.75%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.042'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.68%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.714'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.79%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1151'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.67%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.974'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.11%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3890'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.19%  '
